# ---------------------------------------------------------------
# Quarterly financials update: add two new fiscal-quarter columns
# (2018-09-27 and 2018-12-31), inserted before the existing data,
# and correct several previously-reported historical figures.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank columns at D:E. Excel shifts the existing
#    D:K data (and its formatting) to F:M automatically.
$ws.Columns("D:E").Insert()

# 2) For every row that carries data in columns F:M, copy the cell
#    format from column F into the two new D:E cells (so the new
#    columns inherit the correct date/number style), then write the
#    two new-quarter values (row data transcribed from the source).
$newColData = @{
    7 = @(43465, 43370)
    8 = @(1835300, 1813700)
    9 = @(1534600, 1543100)
    10 = @(300700, 270600)
    11 = @($null, $null)
    12 = @(11200, 10800)
    13 = @(0, 0)
    14 = @(-10000, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(1591700, 1591200)
    18 = @(243600, 222500)
    19 = @($null, $null)
    20 = @(-6200, 7400)
    21 = @(296700, 287800)
    22 = @(19700, 24200)
    23 = @(217700, 205700)
    24 = @(45500, 36900)
    25 = @(0, 0)
    26 = @(172200, 168800)
    27 = @(172500, 168700)
    28 = @(0, 0)
    29 = @(5400, "__STR__10")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(6200, -7400)
    33 = @(177900, 168700)
    34 = @(0, 0)
    35 = @(177900, 168700)
    38 = @(43465, 43370)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(773600, 683400)
    42 = @(0, 0)
    43 = @(1014500, 1198300)
    44 = @(1012600, 931000)
    45 = @(48600, 35700)
    46 = @(2849300, 2848400)
    47 = @(54100, 53500)
    48 = @(2167600, 2123000)
    49 = @(3800, 3900)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(611100, 653300)
    53 = @(0, 0)
    54 = @(5685900, 5682100)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(902600, 882400)
    58 = @(31400, 25000)
    59 = @(648100, 904800)
    60 = @(1582100, 1812200)
    61 = @(1864000, 1869700)
    62 = @(1001700, 883800)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(4448300, 4566200)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(2713200, 2548700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(1237600, 1115900)
    77 = @(0, 0)
    80 = @(43465, 43370)
    81 = @(177900, 168700)
    82 = @($null, $null)
    83 = @(59300, 57900)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(202500, 170200)
    90 = @($null, $null)
    91 = @(-100300, -61500)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-99700, -59000)
    95 = @($null, $null)
    96 = @(-12600, -12600)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-12400, -118100)
    101 = @(-200, 0)
    102 = @(90200, -6900)
}

foreach ($r in $newColData.Keys) {
    $pair = $newColData[$r]
    $ws.Range("F$r").Copy() | Out-Null
    $ws.Range("D${r}:E$r").PasteSpecial(-4122)
    if ($pair[0] -ne $null) {
        if ($pair[0] -is [string] -and $pair[0].StartsWith("__STR__")) {
            $ws.Range("D$r").Value = "NA"
        } else {
            $ws.Range("D$r").Value = $pair[0]
        }
    }
    if ($pair[1] -ne $null) {
        if ($pair[1] -is [string] -and $pair[1].StartsWith("__STR__")) {
            $ws.Range("E$r").Value = "NA"
        } else {
            $ws.Range("E$r").Value = $pair[1]
        }
    }
}

# 3) Apply corrections to a handful of previously-reported values
#    that changed as part of this data refresh (now living in the
#    shifted F:M columns).
$ws.Range("I9").Value = 1486600
$ws.Range("I10").Value = 261600
$ws.Range("I17").Value = 1545900
$ws.Range("I18").Value = 202300
$ws.Range("I20").Value = 11000
$ws.Range("F27").Value = 145100
$ws.Range("I32").Value = -11000
$ws.Range("F33").Value = 145100
$ws.Range("F35").Value = 145100
$ws.Range("F81").Value = 145100
$ws.Range("I91").Value = -50600
$ws.Range("J91").Value = -47500
$ws.Range("H94").Value = -134500
$ws.Range("I94").Value = -50400
$ws.Range("H100").Value = -121200
$ws.Range("I100").Value = -213600
$ws.Range("H102").Value = -306600
$ws.Range("I102").Value = 28600

$ws.Columns("A:M").AutoFit() | Out-Null
